# Fruta / hortaliza, semanal
# Insert two new daily price rows above the existing row 121 (the data block is
# sorted with the newest dates first), shifting the rest of the table down by 2
# rows. The two new rows are populated by duplicating what becomes rows 123/124
# (a copy of the former rows 121/122) and then editing the fields that changed
# for the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 121 - this pushes the former rows 121-212
# down to become rows 123-214.
$ws.Rows("121:122").Insert()

# Seed the two new rows with a copy of the rows now sitting at 123:124 (which
# are exactly the former rows 121:122), then overwrite the cells that actually
# differ for the new entries.
$ws.Range("A123:T124").Copy()
$ws.Range("A121").PasteSpecial()

# New row 121: Ciruela Angeleno, Especial, origin Provincia de Curicó
$ws.Range("D121").Value = 45033
$ws.Range("K121").Value = "Angeleno"
$ws.Range("L121").Value = "Especial"
$ws.Range("M121").Value = 260
$ws.Range("R121").Value = "Provincia de Curicó"

# New row 122: Ciruela Angeleno, Primera, updated price figures
$ws.Range("D122").Value = 45033
$ws.Range("K122").Value = "Angeleno"
$ws.Range("M122").Value = 300
$ws.Range("N122").Value = 10000
$ws.Range("O122").Value = 10000
$ws.Range("P122").Value = 10000
$ws.Range("S122").Value = 556
